$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Append:
#   - one blank paragraph
#   - one blank paragraph with no bullet / no indent
#   - one "Github Repo Of My Code: <url>" paragraph, also with no bullet / no indent
$tr.InsertAfter("`r`r`rGithub Repo Of My Code: https://github.com/Parvat-R/ImageFinder")

$blankPara = $sh.TextFrame.TextRange.Paragraphs(5,1)
$blankPara.ParagraphFormat.Bullet.Visible = 0

$repoPara = $sh.TextFrame.TextRange.Paragraphs(6,1)
$repoPara.ParagraphFormat.Bullet.Visible = 0

# Turn just the URL portion of the last paragraph into a hyperlink
$urlRange = $repoPara.Characters(25, 39)
$urlRange.ActionSettings(1).Hyperlink.Address = "https://github.com/Parvat-R/ImageFinder"
